$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1902.3611
$ws.Range("J17").Value = 1902.3611
$ws.Range("L17").Value = 5707.0833
$ws.Range("N17").Value = -6043.0833
$ws.Range("H69").Value = 7259.25
$ws.Range("I69").Value = 4350.3335
$ws.Range("K69").Value = 13051.0005
$ws.Range("M69").Value = -12177.0005
$ws.Range("H72").Value = 7259.25
$ws.Range("I72").Value = 4350.3335
$ws.Range("K72").Value = 39153.0015
$ws.Range("M72").Value = -34785.0015
$ws.Range("H74").Value = 7836.136
$ws.Range("I74").Value = 4666.6665
$ws.Range("K74").Value = 4666.6665
$ws.Range("M74").Value = -3730.6665
$ws.Range("H77").Value = 7836.136
$ws.Range("I77").Value = 4666.6665
$ws.Range("K77").Value = 23333.3325
$ws.Range("M77").Value = -18653.3325
$ws.Range("H88").Value = 2520.4546
$ws.Range("I88").Value = 3094
$ws.Range("J88").Value = 2192.7144
$ws.Range("K88").Value = 3094
$ws.Range("L88").Value = 2192.7144
$ws.Range("M88").Value = -2688
$ws.Range("N88").Value = -3004.7144
$ws.Range("H91").Value = 2520.4546
$ws.Range("I91").Value = 3094
$ws.Range("J91").Value = 2192.7144
$ws.Range("K91").Value = 3094
$ws.Range("L91").Value = 2192.7144
$ws.Range("M91").Value = -1690
$ws.Range("N91").Value = -5000.7144
$ws.Range("H132").Value = 2482.5518
$ws.Range("I132").Value = 2428.3928
$ws.Range("K132").Value = 7285.178400000001
$ws.Range("M132").Value = -4755.178400000001
$ws.Range("H135").Value = 1272.3784
$ws.Range("I135").Value = 1260.2812
$ws.Range("K135").Value = 11342.5308
$ws.Range("M135").Value = -8807.530799999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3739.9487
$ws.Range("I61").Value = 2336.074
$ws.Range("K61").Value = 2336.074
$ws.Range("M61").Value = -2124.074
$ws.Range("H74").Value = 2193.5854
$ws.Range("I74").Value = 1460.4857
$ws.Range("K74").Value = 1460.4857
$ws.Range("M74").Value = -586.4857
$ws.Range("H77").Value = 2193.5854
$ws.Range("I77").Value = 1460.4857
$ws.Range("K77").Value = 7302.4285
$ws.Range("M77").Value = -2934.4285
$ws.Range("H96").Value = 24171.5
$ws.Range("J96").Value = 24171.5
$ws.Range("L96").Value = 24171.5
$ws.Range("N96").Value = -29663.5
$ws.Range("H122").Value = 4253.7144
$ws.Range("I122").Value = 3935.8
$ws.Range("K122").Value = 11807.4
$ws.Range("M122").Value = -9357.400000000001
$ws.Range("H136").Value = 3739.9487
$ws.Range("I136").Value = 2336.074
$ws.Range("K136").Value = 7008.222
$ws.Range("M136").Value = -4458.222
$ws.Range("H138").Value = 99249.75
$ws.Range("J138").Value = 99249.75
$ws.Range("L138").Value = 99249.75
$ws.Range("N138").Value = -109529.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 60966
$ws.Range("J2").Value = 60966
$ws.Range("L2").Value = 60966
$ws.Range("N2").Value = -61192
$ws.Range("H100").Value = 35000
$ws.Range("J100").Value = 35000
$ws.Range("L100").Value = 35000
$ws.Range("N100").Value = -37164
$ws.Range("H103").Value = 149999
$ws.Range("J103").Value = 149999
$ws.Range("L103").Value = 149999
$ws.Range("N103").Value = -152343
$ws.Range("H105").Value = 1765.4166
$ws.Range("I105").Value = 1643.85
$ws.Range("K105").Value = 1643.85
$ws.Range("M105").Value = 103.1500000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 303.875
$ws.Range("I7").Value = 58.333332
$ws.Range("J7").Value = 619.5714
$ws.Range("K7").Value = 58.333332
$ws.Range("L7").Value = 619.5714
$ws.Range("M7").Value = 54.666668
$ws.Range("N7").Value = -845.5714
$ws.Range("H21").Value = 5000
$ws.Range("I21").Value = 5000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 5000
$ws.Range("L21").Value = 0
$ws.Range("N21").Value = -4765
$ws.Range("M21").ClearContents()
$ws.Range("H69").Value = 50239.11
$ws.Range("I69").Value = 51634.4
$ws.Range("J69").Value = 48495
$ws.Range("K69").Value = 51634.4
$ws.Range("L69").Value = 48495
$ws.Range("M69").Value = -50885.4
$ws.Range("N69").Value = -49993
$ws.Range("H72").Value = 50239.11
$ws.Range("I72").Value = 51634.4
$ws.Range("J72").Value = 48495
$ws.Range("K72").Value = 154903.2
$ws.Range("L72").Value = 145485
$ws.Range("M72").Value = -151159.2
$ws.Range("N72").Value = -152973
$ws.Range("H97").Value = 59777.5
$ws.Range("J97").Value = 59777.5
$ws.Range("L97").Value = 59777.5
$ws.Range("N97").Value = -61759.5
$ws.Range("H134").Value = 1860.3667
$ws.Range("I134").Value = 1840.8462
$ws.Range("K134").Value = 5522.5386
$ws.Range("M134").Value = -2987.5386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 264.7
$ws.Range("I11").Value = 130.875
$ws.Range("J11").Value = 800
$ws.Range("K11").Value = 392.625
$ws.Range("L11").Value = 2400
$ws.Range("M11").Value = -252.625
$ws.Range("N11").Value = -2680
$ws.Range("H63").Value = 12511.223
$ws.Range("I63").Value = 8000
$ws.Range("J63").Value = 13075.125
$ws.Range("K63").Value = 24000
$ws.Range("L63").Value = 39225.375
$ws.Range("M63").Value = -23251
$ws.Range("N63").Value = -40723.375
$ws.Range("H66").Value = 12511.223
$ws.Range("I66").Value = 8000
$ws.Range("J66").Value = 13075.125
$ws.Range("K66").Value = 72000
$ws.Range("L66").Value = 117676.125
$ws.Range("M66").Value = -68256
$ws.Range("N66").Value = -125164.125
$ws.Range("H69").Value = 4170.7144
$ws.Range("J69").Value = 6548.75
$ws.Range("L69").Value = 19646.25
$ws.Range("N69").Value = -21268.25
$ws.Range("H72").Value = 4170.7144
$ws.Range("J72").Value = 6548.75
$ws.Range("L72").Value = 58938.75
$ws.Range("N72").Value = -67050.75
$ws.Range("H131").Value = 1323.7646
$ws.Range("J131").Value = 2162.25
$ws.Range("L131").Value = 6486.75
$ws.Range("N131").Value = -16566.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 6295.2
$ws.Range("I17").Value = 5998
$ws.Range("J17").Value = 6493.3335
$ws.Range("K17").Value = 5998
$ws.Range("L17").Value = 6493.3335
$ws.Range("M17").Value = -5830
$ws.Range("N17").Value = -6829.3335
$ws.Range("H39").Value = 23129.5
$ws.Range("J39").Value = 23129.5
$ws.Range("L39").Value = 23129.5
$ws.Range("N39").Value = -24193.5
$ws.Range("H80").Value = 4681.9443
$ws.Range("I80").Value = 4297.5835
$ws.Range("J80").Value = 5450.6665
$ws.Range("K80").Value = 4297.5835
$ws.Range("L80").Value = 5450.6665
$ws.Range("M80").Value = -3299.5835
$ws.Range("N80").Value = -7446.6665
$ws.Range("H83").Value = 4681.9443
$ws.Range("I83").Value = 4297.5835
$ws.Range("J83").Value = 5450.6665
$ws.Range("K83").Value = 21487.9175
$ws.Range("L83").Value = 27253.3325
$ws.Range("M83").Value = -16495.9175
$ws.Range("N83").Value = -37237.3325
$ws.Range("H86").Value = 90143
$ws.Range("J86").Value = 90143
$ws.Range("L86").Value = 90143
$ws.Range("N86").Value = -92515
$ws.Range("H89").Value = 90143
$ws.Range("J89").Value = 90143
$ws.Range("L89").Value = 270429
$ws.Range("N89").Value = -282285
$ws.Range("H92").Value = 56077.777
$ws.Range("I92").Value = 4000
$ws.Range("J92").Value = 62587.5
$ws.Range("K92").Value = 4000
$ws.Range("L92").Value = 62587.5
$ws.Range("M92").Value = -2128
$ws.Range("N92").Value = -66331.5
$ws.Range("H97").Value = 1724.9395
$ws.Range("I97").Value = 641.0476
$ws.Range("J97").Value = 3621.75
$ws.Range("K97").Value = 641.0476
$ws.Range("L97").Value = 3621.75
$ws.Range("M97").Value = -145.0476
$ws.Range("N97").Value = -4613.75
$ws.Range("H102").Value = 2195.125
$ws.Range("H132").Value = 1034.625
$ws.Range("I132").Value = 1006.0769
$ws.Range("K132").Value = 3018.2307
$ws.Range("M132").Value = -488.2307000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 101699.1
$ws.Range("I22").Value = 167333.5
$ws.Range("J22").Value = 3247.5
$ws.Range("K22").Value = 167333.5
$ws.Range("L22").Value = 3247.5
$ws.Range("M22").Value = -167038.5
$ws.Range("N22").Value = -3837.5
$ws.Range("H27").Value = 101699.1
$ws.Range("I27").Value = 167333.5
$ws.Range("J27").Value = 3247.5
$ws.Range("K27").Value = 167333.5
$ws.Range("L27").Value = 3247.5
$ws.Range("M27").Value = -167226.5
$ws.Range("N27").Value = -3461.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 78493.875
$ws.Range("I75").Value = 77847.5
$ws.Range("K75").Value = 77847.5
$ws.Range("M75").Value = -76911.5
$ws.Range("H78").Value = 78493.875
$ws.Range("I78").Value = 77847.5
$ws.Range("K78").Value = 233542.5
$ws.Range("M78").Value = -228862.5
